$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prices (column D) and links (column F) for each existing row
$ws.Range("D2").Value = "48,35TL "
$ws.Range("D3").Value = "11,73TL"
$ws.Range("D4").Value = "36,45TL"
$ws.Range("D5").Value = "4,69TL"
$ws.Range("D6").Value = "25TL"
$ws.Range("D7").Value = "2,4TL"
$ws.Range("D8").Value = "14,65TL"

$ws.Range("F1").Value = "Link"
$ws.Range("F2").Value = "https://ozdisan.com/pasif-komponentler/kondansatorler/aluminyum-kondansatorler/PKL5-400V221MN400"
$ws.Range("F3").Value = "https://www.direnc.net/35a-1000v-kopru-diyot"
$ws.Range("F4").Value = "https://www.direnc.net/ixgh24n60c4d1-rohs-24a600v-to247ad-igbtdiode"
$ws.Range("F5").Value = "https://www.direnc.net/tl494--025a-switching-controller-300khz-switching-freq-max"
$ws.Range("F6").Value = "https://ozdisan.com/guc-yari-iletkenleri/diyotlar-modul-diyotlar-ve-dogrultucular/genel-amacli-diyotlar/DHG30I600PA"
$ws.Range("F7").Value = "https://ozdisan.com/pasif-komponentler/kondansatorler/film-kondansatorler/C322J104J60A605"
$ws.Range("F8").Value = "https://ozdisan.com/pasif-komponentler/direncler/tas-direncler/PRW05WJW10KB00"

$ws.Columns.Item(6).ColumnWidth = 128.140625

# Remove the old hyperlink in B17 but keep its (hyperlink) style
$ws.Range("B17").Hyperlinks.Delete()
$ws.Range("B17").Value = $null

$ws.Range("D11").Select()
